# Insert a new data row at row 221 (pushes existing rows 221..261 down to 222..262)
# and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(221).Insert()

$ws.Range("A221").Value = 10
$ws.Range("B221").Value = "Vega Modelo de Temuco"
$ws.Range("C221").Value = "La Araucanía"
$ws.Range("D221").Value = 44798
$ws.Range("E221").Value = 9
$ws.Range("F221").Value = 100112052
$ws.Range("G221").Value = "Albahaca"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 50
$ws.Range("K221").Value = 6000
$ws.Range("L221").Value = 6000
$ws.Range("M221").Value = 6000
$ws.Range("N221").Value = "$/paquete"
$ws.Range("O221").Value = "Región de Arica y Parinacota"
$ws.Range("P221").Value = 6000
$ws.Range("Q221").Value = 1
$ws.Range("R221").Value = "Hortaliza"
